$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F2: model_name  fsdfsd -> fds
$ws.Range("F2").Value = "fds"

# G2: kolvo  432 -> 435 (keep as text, matches original inlineStr type)
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "435"

# H2: call_me  gfdsfsd -> 432432 (keep as text)
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "432432"

# J2: power  Подключиться к пулу со скидкой ✅ -> Хочу купить оборудование 🔥
$ws.Range("J2").Value = "Хочу купить оборудование 🔥"

# K2: currency  USA -> USD $
$ws.Range("K2").Value = "USD $"

# M2: cost_electricity  43 -> 4343 (keep as text)
$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = "4343"

# N2: hash  3 -> 43 (keep as text)
$ws.Range("N2").NumberFormat = "@"
$ws.Range("N2").Value = "43"

# O2: potreb  2 -> 43 (keep as text)
$ws.Range("O2").NumberFormat = "@"
$ws.Range("O2").Value = "43"

# P2: komm  /password -> 43 (keep as text)
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "43"

# Q2: promo  yes -> да
$ws.Range("Q2").Value = "да"

# R2: date  2023-06-12 23:36:23.261064 -> 2023-06-13 10:54:32.282695
$ws.Range("R2").Value = "2023-06-13 10:54:32.282695"
